# Add a new "2022-Q3" sheet right after "总计", with its fund-holding data,
# and record the new quarter's summary figures on the "总计" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update the summary sheet ("总计"): insert a new row for 2022-Q3 and
#    shift the existing quarter rows down by one.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Push all existing data rows (old rows 2..8) down by one row.
$summary.Rows.Item(2).Insert()

# Fill in the brand-new 2022-Q3 row.
$summary.Cells.Item(2, 1).Value = 1
$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 3).Value = 6
$summary.Cells.Item(2, 4).Value = 1.1

# Fix up the running index in column A for the rows that got shifted down
# (they must read 1..7 now instead of 0..6).
$summary.Cells.Item(3, 1).Value = 1
$summary.Cells.Item(4, 1).Value = 2
$summary.Cells.Item(5, 1).Value = 3
$summary.Cells.Item(6, 1).Value = 4
$summary.Cells.Item(7, 1).Value = 5
$summary.Cells.Item(8, 1).Value = 6
$summary.Cells.Item(9, 1).Value = 7

# ---------------------------------------------------------------------
# 2. Create the new "2022-Q3" worksheet (placed right after "总计") and
#    populate it with the quarter's fund holdings table.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $summary)
$q3.Name = "2022-Q3"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 0; $col -lt $headers.Length; $col++) {
    $q3.Cells.Item(1, $col + 2).Value = $headers[$col]
}

# Columns: index(A,n) code(B,text) name(C,text) size(D,text) position(E,text)
#          ratio(F,text) marketValue(G,text) rank(H,n)
$rows = @(
    @(0, "161219", "国投瑞银新兴产业混合（LOF）",             "6.18", "79.94", "4.54", "0.2806", 5),
    @(1, "161232", "国投瑞银瑞盛灵活配置混合",                 "4.18", "94.55", "6.71", "0.2805", 3),
    @(2, "000663", "国投瑞银美丽中国灵活配置混合",             "3.85", "93.40", "6.91", "0.2660", 2),
    @(3, "161225", "国投瑞银瑞盈灵活配置混合（LOF）",           "2.52", "94.35", "7.39", "0.1862", 2),
    @(4, "161233", "国投瑞银瑞泰多策略灵活配置混合（LOF）A",    "4.96", "28.33", "1.56", "0.0774", 5),
    @(5, "011618", "国投瑞银瑞泰多策略灵活配置混合（LOF）C",    "0.80", "28.33", "1.56", "0.0125", 5)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]
    $rowNum = $i + 2
    $q3.Cells.Item($rowNum, 1).Value = $r[0]
    $q3.Cells.Item($rowNum, 2).Value = "'" + $r[1]
    $q3.Cells.Item($rowNum, 3).Value = $r[2]
    $q3.Cells.Item($rowNum, 4).Value = "'" + $r[3]
    $q3.Cells.Item($rowNum, 5).Value = "'" + $r[4]
    $q3.Cells.Item($rowNum, 6).Value = "'" + $r[5]
    $q3.Cells.Item($rowNum, 7).Value = "'" + $r[6]
    $q3.Cells.Item($rowNum, 8).Value = $r[7]
}

Write-Output "2022-Q3 sheet added and 总计 sheet updated."
